$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text/value updates (non-ambiguous strings; Excel keeps these as text automatically)
$ws.Range("D2").Value = "24.918.45"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.707.39"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +5.41%  "
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +5.34%  "
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("E15").Value = "  +4.91%  "
$ws.Range("D16").Value = "1.709.46"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("E21").Value = "  +5.70%  "
$ws.Range("E22").Value = "  +4.70%  "
$ws.Range("E23").Value = "  +8.04%  "
$ws.Range("D24").Value = "24.926.56"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +6.05%  "
$ws.Range("E27").Value = "  +5.05%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  +3.97%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.894.83"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E31").Value = "  +29.07%  "
$ws.Range("E32").Value = "  +8.57%  "
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("E34").Value = "  +13.28%  "
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("E41").Value = "  +7.37%  "
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("E43").Value = "  +6.70%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("E45").Value = "  +6.77%  "
$ws.Range("E46").Value = "  +5.12%  "
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("E48").Value = "  +5.27%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("E51").Value = "  +6.30%  "

# Numeric-looking price strings must be forced to text without altering cell style.
# Trick: build the literal string via a formula in a scratch cell, copy it, then
# paste-special VALUES ONLY into the target cell (keeps style untouched), then wipe the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1.005"""
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Formula = "=""312.98"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Formula = "=""0.9989"""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Formula = "=""0.3744"""
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Formula = "=""49.35"""
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Formula = "=""0.3440"""
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Formula = "=""1.224"""
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Formula = "=""0.07536"""
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Formula = "=""1.000"""
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Formula = "=""21.19"""
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Formula = "=""6.324"""
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Formula = "=""7.068"""
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Formula = "=""0.00001130"""
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Formula = "=""0.06726"""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$scratch.Formula = "=""0.9978"""
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Formula = "=""84.01"""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Formula = "=""6.380"""
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Formula = "=""2.442"""
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Formula = "=""2.800"""
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Formula = "=""20.39"""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Formula = "=""150.14"""
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Formula = "=""132.78"""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Formula = "=""1.250"""
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Formula = "=""6.805"""
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Formula = "=""13.82"""
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Formula = "=""0.08802"""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Formula = "=""1.772"""
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Formula = "=""5.627"""
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Formula = "=""0.06663"""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Formula = "=""9.135"""
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Formula = "=""0.02418"""
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Formula = "=""0.2242"""
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Formula = "=""1.271"""
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Formula = "=""0.6501"""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Formula = "=""0.9983"""
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Formula = "=""0.6178"""
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Formula = "=""3.839"""
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Formula = "=""2.120"""
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Formula = "=""129.15"""
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Formula = "=""0.07336"""
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
